$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new date
$ws.Name = "6-8-24"

# Clear all existing content so that the shared-strings table gets rebuilt from
# scratch in exactly the order we (re)introduce values below.
$ws.Range("A1:B11").ClearContents()

# --- Row 1 (headers) ---
$ws.Range("A1").Value = "Games"
$ws.Range("B1").Value = "RFPred"
$ws.Range("C1").Value = "NRFI"
$ws.Range("D1").Value = "Correct"
$ws.Range("E1").Value = "Total"
$ws.Range("G1").Value = "Percent"

# --- Data rows (Games / RFPred) ---
$ws.Range("A2").Value = "('KC', 'SEA')"
$ws.Range("B2").Value = 0.934

$ws.Range("A3").Value = "('MIN', 'PIT')"
$ws.Range("B3").Value = 0.834

$ws.Range("A4").Value = "('LAD', 'NYY')"
$ws.Range("B4").Value = 0.792

$ws.Range("A5").Value = "('NYM', 'PHI')"
$ws.Range("B5").Value = 0.786

$ws.Range("A6").Value = "('DET', 'MIL')"
$ws.Range("B6").Value = 0.775

$ws.Range("A7").Value = "('SF', 'TEX')"
$ws.Range("B7").Value = 0.724

$ws.Range("A8").Value = "('BAL', 'TB')"
$ws.Range("B8").Value = 0.675

$ws.Range("A9").Value = "('ATL', 'WSH')"
$ws.Range("B9").Value = 0.555

$ws.Range("A10").Value = "('BOS', 'CWS')"
$ws.Range("B10").Value = 0.545

$ws.Range("A11").Value = "('CHC', 'CIN')"
$ws.Range("B11").Value = 0.513

$ws.Range("A12").Value = "('AZ', 'SD')"
$ws.Range("B12").Value = 0.484

$ws.Range("A13").Value = "('HOU', 'LAA')"
$ws.Range("B13").Value = 0.477

$ws.Range("A14").Value = "('OAK', 'TOR')"
$ws.Range("B14").Value = 0.328

$ws.Range("A15").Value = "('COL', 'STL')"
$ws.Range("B15").Value = 0.32

$ws.Range("A16").Value = "('CLE', 'MIA')"
$ws.Range("B16").Value = 0.049

# Restore the selection to match the post-edit state
$ws.Range("E9").Select()
